# Re-upload of the roster data: the player/position/team table is
# reshuffled into a new canonical order (matching the order the shared
# strings now appear in), even though the set of rows (and the sheet's
# row/column layout) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @('Jalen Brunson', 'Devin Booker', 'Desmond Bane', 'Jrue Holiday', 'Shaedon Sharpe', 'Jeremy Sochan', 'Devin Vassell', 'Draymond Green', 'Yves Missi', 'D''Angelo Russell', 'Alperen Sengün', 'Trae Young', 'Coby White', 'Norman Powell', 'Walker Kessler', 'Immanuel Quickley', 'Kawhi Leonard', 'LeBron James')
$positions = @('PG', 'PG,SG', 'SG,SF', 'PG,SG', 'SG,SF', 'SF,PF', 'SG,SF', 'PF,C', 'C', 'PG', 'C', 'PG', 'PG,SG', 'SG,SF', 'C', 'PG,SG', 'SG,SF,PF', 'SF,PF')
$teams = @('New York Knicks', 'Phoenix Suns', 'Memphis Grizzlies', 'Boston Celtics', 'Portland Trail Blazers', 'San Antonio Spurs', 'San Antonio Spurs', 'Golden State Warriors', 'New Orleans Pelicans', 'Los Angeles Lakers', 'Houston Rockets', 'Atlanta Hawks', 'Chicago Bulls', 'LA Clippers', 'Utah Jazz', 'Toronto Raptors', 'LA Clippers', 'Los Angeles Lakers')

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $positions[$i]
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
